$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 516, pushing rows 516..566 down to 517..567
$ws.Rows(516).Insert()

# Fill in the new dataset entry: "Thin Anti-Establishment Supply Dataset"
$ws.Range("A516").Value = "Thin Anti-Establishment Supply Dataset"
$ws.Range("B516").Value = "parties and politicians"
$ws.Range("C516").Value = "https://doi.org/10.1177/13540688221080536"
$ws.Range("D516").Value = "anti-establishment politics, election campaign, populism"

$ws.Range("F516").Value = 0
$ws.Range("G516").Value = 0
$ws.Range("H516").Value = 1
$ws.Range("I516").Value = 0
$ws.Range("J516").Value = 1

$ws.Range("K516").Value = 2010
$ws.Range("L516").Value = 2019

$ws.Range("M516").Value = "online"
$ws.Range("N516").Value = "no"
$ws.Range("O516").Value = "CC-BY-4.0"

$ws.Range("Z516").Value = "10.1177/13540688221080536"
$ws.Range("AA516").Value = "10.17605/OSF.IO/F23HM"

$ws.Range("AB516").Value = 20240114

# Turn the link column into a real hyperlink, then restore the sheet's
# standard hyperlink cell style (Hyperlinks.Add applies its own style).
$ws.Hyperlinks.Add($ws.Range("C516"), "https://doi.org/10.1177/13540688221080536")
$ws.Range("C516").Style = "Hyperlink"

# Move the view/selection the way the author's session ended up
$ws.Application.ActiveWindow.ScrollRow = 542
$ws.Range("A563").Select()
